$d = $word.ActiveDocument

# 1. Fix the title paragraph's paragraph-mark formatting (first paragraph).
#    Remove eastAsia font / color / underline / lang, change size 32 -> 28 (+ add szCs 28).
$titlePara = $d.Paragraphs(1)
$titleMark = $titlePara.Range.ParagraphFormat
$markFont = $titlePara.Range.Font
# The paragraph mark's own character formatting is accessed via the paragraph's
# Range that covers just the end-of-paragraph mark.
$markRange = $titlePara.Range
$markRange.SetRange($markRange.End - 1, $markRange.End)
$markRange.Font.Name = "Arial"
$markRange.Font.NameFarEast = ""
$markRange.Font.NameAscii = "Arial"
$markRange.Font.Color = 0
$markRange.Font.Size = 14
$markRange.Font.Underline = 0
$markRange.Font.Bold = 0

# 2. Remove the "(6 seconds)" and "(4 seconds)" parenthetical remarks.
$d.Content.Find.Execute("parallelisable part of the program takes up 60% of the program (6 seconds), the parallelisable part takes up 40% of the program (4 seconds) how long will the program take to run on 4 cores?", $true, $false, $false, $false, $false, $true, 1, $false, "parallelisable part of the program takes up 60% of the program, the parallelisable part takes up 40% of the program how long will the program take to run on 4 cores?", 2)

# 3. Remove the two trailing empty paragraphs at the end of the document.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$secondLastPara = $d.Paragraphs($count - 1)
$removeRange = $d.Range($secondLastPara.Range.Start, $lastPara.Range.End)
$removeRange.Delete()
